{"js": "// Update the date paragraph and the 25 division problems in the table,\n// following the document order of the text runs (there are duplicate\n// problem strings, so replacement must be positional, not text-matched).\n\nconst body = context.document.body;\n\n// 1) Date heading paragraph (first paragraph of the body).\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst datePara = paras.items[0];\ndatePara.insertText(\"2025-11-08 Saturday\", Word.InsertLocation.replace);\n\n// 2) Division problems inside the table \u2014 only rows 0,4,8,12,16 (of 20)\n// contain text; the other rows are blank spacer rows. Each populated row\n// has 5 cells (columns 0-4), giving 25 problems in document order.\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newValues = [\n  [\"51\u00f79=\", \"86\u00f75=\", \"52\u00f73=\", \"53\u00f77=\", \"81\u00f76=\"],\n  [\"51\u00f72=\", \"72\u00f74=\", \"45\u00f74=\", \"57\u00f78=\", \"33\u00f77=\"],\n  [\"38\u00f76=\", \"82\u00f78=\", \"26\u00f76=\", \"41\u00f79=\", \"30\u00f74=\"],\n  [\"49\u00f73=\", \"94\u00f73=\", \"51\u00f72=\", \"76\u00f79=\", \"93\u00f75=\"],\n  [\"39\u00f79=\", \"48\u00f72=\", \"86\u00f76=\", \"47\u00f74=\", \"65\u00f77=\"],\n];\n\nconst contentRows = [0, 4, 8, 12, 16];\n\nfor (let i = 0; i < contentRows.length; i++) {\n  const rowIndex = contentRows[i];\n  for (let col = 0; col < 5; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = newValues[i][col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 division problems in the table.\n# Word COM is 1-indexed. Text is addressed positionally (by paragraph /\n# row / column) because several problem strings repeat in the document,\n# so a blind Find-and-ReplaceAll would touch the wrong occurrences.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading \u2014 first paragraph of the document.\n$d.Paragraphs.Item(1).Range.Text = \"2025-11-08 Saturday\"\n\n# 2) Division problems inside the table. Only rows 1,5,9,13,17 (of 20)\n# hold text; the rest are blank spacer rows. Each populated row has 5\n# cells, for 25 problems total, read/written in document order.\n$t = $d.Tables.Item(1)\n\n$newValues = @(\n    @(\"51\u00f79=\", \"86\u00f75=\", \"52\u00f73=\", \"53\u00f77=\", \"81\u00f76=\"),\n    @(\"51\u00f72=\", \"72\u00f74=\", \"45\u00f74=\", \"57\u00f78=\", \"33\u00f77=\"),\n    @(\"38\u00f76=\", \"82\u00f78=\", \"26\u00f76=\", \"41\u00f79=\", \"30\u00f74=\"),\n    @(\"49\u00f73=\", \"94\u00f73=\", \"51\u00f72=\", \"76\u00f79=\", \"93\u00f75=\"),\n    @(\"39\u00f79=\", \"48\u00f72=\", \"86\u00f76=\", \"47\u00f74=\", \"65\u00f77=\")\n)\n\n$contentRows = @(1, 5, 9, 13, 17)\n\nfor ($i = 0; $i -lt $contentRows.Length; $i++) {\n    $rowIndex = $contentRows[$i]\n    for ($col = 1; $col -le 5; $col++) {\n        $cell = $t.Cell($rowIndex, $col)\n        $cell.Range.Text = $newValues[$i][$col - 1]\n    }\n}\n"}
